$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.536.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.72%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.300.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.67%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.630'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.97'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.67%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.419'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0922'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.68%  '

$ws.Range('E11').Value = '  +0.95%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.628.55'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.10%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.95'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.64%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.79%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.817'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.296.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.17%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.411.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.55%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0934'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.51%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.27%  '

$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.60'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.59'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.21%  '

$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.62%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.15%  '

$ws.Range('E25').Value = '  +1.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.53%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.143'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.77%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.68%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.122'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.69%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.50%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0661'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.99%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.44%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.37%  '

$ws.Range('E39').Value = '  +5.45%  '

$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.82'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.09%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000223'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -11.59%  '

$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.72%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.18%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0970'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.49%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.481.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.99%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.52%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.81%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.28'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.64%  '

$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.71%  '
